$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear row 7 columns A:H (no longer present in the re-run results) ---
$ws.Range("A7:H7").Clear()

# Row 1
$ws.Range("A1").Value = 'negative'
$ws.Range("J1").Value = 'positive'

# Row 2
$ws.Range("A2").Value = 'name'
$ws.Range("B2").Value = 'anchor score'
$ws.Range("C2").Value = 'type occurences'
$ws.Range("D2").Value = 'total occurences'
$ws.Range("E2").Value = '+%'
$ws.Range("F2").Value = '-%'
$ws.Range("G2").Value = 'both'
$ws.Range("H2").Value = 'normal'
$ws.Range("J2").Value = 'name'
$ws.Range("K2").Value = 'anchor score'
$ws.Range("L2").Value = 'type occurences'
$ws.Range("M2").Value = 'total occurences'
$ws.Range("N2").Value = '+%'
$ws.Range("O2").Value = '-%'
$ws.Range("P2").Value = 'both'
$ws.Range("Q2").Value = 'normal'

# Row 3
$ws.Range("A3").Value = 'crude'
$ws.Range("B3").Value = 0.8823529411764706
$ws.Range("C3").Value = 30
$ws.Range("D3").Value = 30
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = 4
$ws.Range("J3").Value = 'happy'
$ws.Range("K3").Value = 0.9615384615384616
$ws.Range("L3").Value = 25
$ws.Range("M3").Value = 25
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = 1

# Row 4
$ws.Range("A4").Value = 'crisis'
$ws.Range("B4").Value = 0.6232876712328768
$ws.Range("C4").Value = 182
$ws.Range("D4").Value = 182
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = 110
$ws.Range("J4").Value = 'best'
$ws.Range("K4").Value = 0.9322033898305084
$ws.Range("L4").Value = 55
$ws.Range("M4").Value = 55
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = 4

# Row 5
$ws.Range("A5").Value = 'sc'
$ws.Range("B5").Value = 0.2222222222222222
$ws.Range("C5").Value = 42
$ws.Range("D5").Value = 42
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 147
$ws.Range("J5").Value = 'interesting'
$ws.Range("K5").Value = 0.9090909090909091
$ws.Range("L5").Value = 30
$ws.Range("M5").Value = 30
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = 3

# Row 6
$ws.Range("A6").Value = 'panic'
$ws.Range("B6").Value = 0.2131782945736434
$ws.Range("C6").Value = 110
$ws.Range("D6").Value = 110
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = 406
$ws.Range("J6").Value = 'love'
$ws.Range("K6").Value = 0.8695652173913043
$ws.Range("L6").Value = 40
$ws.Range("M6").Value = 40
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = $false
$ws.Range("Q6").Value = 6

# Row 7
$ws.Range("J7").Value = 'great'
$ws.Range("K7").Value = 0.8482142857142857
$ws.Range("L7").Value = 95
$ws.Range("M7").Value = 95
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = 17

# Row 8
$ws.Range("J8").Value = 'thanks'
$ws.Range("K8").Value = 0.8414634146341463
$ws.Range("L8").Value = 69
$ws.Range("M8").Value = 69
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = $false
$ws.Range("Q8").Value = 13

# Row 9
$ws.Range("J9").Value = 'thank'
$ws.Range("K9").Value = 0.7890625
$ws.Range("L9").Value = 101
$ws.Range("M9").Value = 101
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = $false
$ws.Range("Q9").Value = 27

# Row 10
$ws.Range("J10").Value = 'free'
$ws.Range("K10").Value = 0.775
$ws.Range("L10").Value = 93
$ws.Range("M10").Value = 93
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = $false
$ws.Range("Q10").Value = 27

# Row 11
$ws.Range("J11").Value = 'positive'
$ws.Range("K11").Value = 0.7413793103448276
$ws.Range("L11").Value = 43
$ws.Range("M11").Value = 43
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = $false
$ws.Range("Q11").Value = 15

# Row 12
$ws.Range("J12").Value = 'safe'
$ws.Range("K12").Value = 0.7253521126760564
$ws.Range("L12").Value = 103
$ws.Range("M12").Value = 103
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = $false
$ws.Range("Q12").Value = 39

# Row 13
$ws.Range("J13").Value = 'special'
$ws.Range("K13").Value = 0.6944444444444444
$ws.Range("L13").Value = 25
$ws.Range("M13").Value = 25
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = $false
$ws.Range("Q13").Value = 11

# Row 14
$ws.Range("J14").Value = 'support'
$ws.Range("K14").Value = 0.6886792452830188
$ws.Range("L14").Value = 73
$ws.Range("M14").Value = 73
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = $false
$ws.Range("Q14").Value = 33

# Row 15
$ws.Range("J15").Value = 'good'
$ws.Range("K15").Value = 0.6875
$ws.Range("L15").Value = 110
$ws.Range("M15").Value = 110
$ws.Range("N15").Value = 1
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = $false
$ws.Range("Q15").Value = 50

# Row 16
$ws.Range("J16").Value = 'relief'
$ws.Range("K16").Value = 0.62
$ws.Range("L16").Value = 31
$ws.Range("M16").Value = 31
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = $false
$ws.Range("Q16").Value = 19

# Row 17
$ws.Range("J17").Value = 'well'
$ws.Range("K17").Value = 0.6063829787234043
$ws.Range("L17").Value = 57
$ws.Range("M17").Value = 57
$ws.Range("N17").Value = 1
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = $false
$ws.Range("Q17").Value = 37

# Row 18
$ws.Range("J18").Value = 'fresh'
$ws.Range("K18").Value = 0.6041666666666666
$ws.Range("L18").Value = 29
$ws.Range("M18").Value = 29
$ws.Range("N18").Value = 1
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = $false
$ws.Range("Q18").Value = 19

# Row 19
$ws.Range("J19").Value = 'heroes'
$ws.Range("K19").Value = 0.5957446808510638
$ws.Range("L19").Value = 28
$ws.Range("M19").Value = 28
$ws.Range("N19").Value = 1
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = $false
$ws.Range("Q19").Value = 19

# Row 20
$ws.Range("J20").Value = 'safety'
$ws.Range("K20").Value = 0.5882352941176471
$ws.Range("L20").Value = 30
$ws.Range("M20").Value = 30
$ws.Range("N20").Value = 1
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = $false
$ws.Range("Q20").Value = 21

# Row 21
$ws.Range("J21").Value = 'better'
$ws.Range("K21").Value = 0.5714285714285714
$ws.Range("L21").Value = 36
$ws.Range("M21").Value = 36
$ws.Range("N21").Value = 1
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = $false
$ws.Range("Q21").Value = 27

# Row 22
$ws.Range("J22").Value = 'hand'
$ws.Range("K22").Value = 0.5169712793733682
$ws.Range("L22").Value = 198
$ws.Range("M22").Value = 198
$ws.Range("N22").Value = 1
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = $false
$ws.Range("Q22").Value = 185

# Row 23
$ws.Range("J23").Value = 'like'
$ws.Range("K23").Value = 0.4823529411764706
$ws.Range("L23").Value = 164
$ws.Range("M23").Value = 164
$ws.Range("N23").Value = 1
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = $false
$ws.Range("Q23").Value = 176

# Row 24
$ws.Range("J24").Value = 'care'
$ws.Range("K24").Value = 0.4719101123595505
$ws.Range("L24").Value = 42
$ws.Range("M24").Value = 42
$ws.Range("N24").Value = 1
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = $false
$ws.Range("Q24").Value = 47

# Row 25
$ws.Range("J25").Value = 'help'
$ws.Range("K25").Value = 0.4406779661016949
$ws.Range("L25").Value = 130
$ws.Range("M25").Value = 130
$ws.Range("N25").Value = 1
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = $false
$ws.Range("Q25").Value = 165

# Row 26
$ws.Range("J26").Value = 'protect'
$ws.Range("K26").Value = 0.4246575342465753
$ws.Range("L26").Value = 31
$ws.Range("M26").Value = 31
$ws.Range("N26").Value = 1
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = $false
$ws.Range("Q26").Value = 42

# Row 27
$ws.Range("J27").Value = 'hope'
$ws.Range("K27").Value = 0.4153846153846154
$ws.Range("L27").Value = 27
$ws.Range("M27").Value = 27
$ws.Range("N27").Value = 1
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = $false
$ws.Range("Q27").Value = 38

# Row 28
$ws.Range("J28").Value = 'sure'
$ws.Range("K28").Value = 0.40625
$ws.Range("L28").Value = 26
$ws.Range("M28").Value = 26
$ws.Range("N28").Value = 1
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = $false
$ws.Range("Q28").Value = 38

# Row 29
$ws.Range("J29").Value = 'increase'
$ws.Range("K29").Value = 0.3717948717948718
$ws.Range("L29").Value = 29
$ws.Range("M29").Value = 29
$ws.Range("N29").Value = 1
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = $false
$ws.Range("Q29").Value = 49

# Row 30
$ws.Range("J30").Value = 'please'
$ws.Range("K30").Value = 0.3682008368200837
$ws.Range("L30").Value = 88
$ws.Range("M30").Value = 88
$ws.Range("N30").Value = 1
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = $false
$ws.Range("Q30").Value = 151

# Row 31
$ws.Range("J31").Value = 'store'
$ws.Range("K31").Value = 0.03691275167785235
$ws.Range("L31").Value = 33
$ws.Range("M31").Value = 33
$ws.Range("N31").Value = 1
$ws.Range("O31").Value = 0
$ws.Range("P31").Value = $false
$ws.Range("Q31").Value = 861

# Row 32
$ws.Range("J32").NumberFormat = "@"
$ws.Range("J32").Value = '19'
$ws.Range("J32").Style = "Normal"
$ws.Range("K32").Value = 0.01635514018691589
$ws.Range("L32").Value = 35
$ws.Range("M32").Value = 38
$ws.Range("N32").Value = 0.92
$ws.Range("O32").Value = 0.07999999999999996
$ws.Range("P32").Value = $true
$ws.Range("Q32").Value = 2105

# Row 33
$ws.Range("J33").Value = 'co'
$ws.Range("K33").Value = 0.01129396579541788
$ws.Range("L33").Value = 35
$ws.Range("M33").Value = 43
$ws.Range("N33").Value = 0.81
$ws.Range("O33").Value = 0.1899999999999999
$ws.Range("P33").Value = $true
$ws.Range("Q33").Value = 3064

# --- Apply the bold + thin-border + center/top style ("style 1") to label cells ---
$r = $ws.Range("A2,B2,C2,D2,E2,F2,G2,H2,J2,K2,L2,M2,N2,O2,P2,Q2,A3,J3,A4,J4")
$r.Font.Bold = $true
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108
$r.VerticalAlignment = -4160

$r = $ws.Range("A5,J5,A6,J6,J7,J8,J9,J10,J11,J12,J13,J14,J15,J16,J17,J18,J19,J20,J21,J22")
$r.Font.Bold = $true
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108
$r.VerticalAlignment = -4160

$r = $ws.Range("J23,J24,J25,J26,J27,J28,J29,J30,J31,J32,J33")
$r.Font.Bold = $true
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108
$r.VerticalAlignment = -4160

Write-Output "done"